$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash a copy of the current "False" text (row 2) into an unused scratch
# cell before it gets overwritten below.
$ws.Range("C2").Copy()
$ws.Range("E1").PasteSpecial(-4163)

# Local config block (rows 2-5) becomes the active one: Usar = True
# (use C6, which already holds the text "True", as the copy source so the
# written cells stay text cells referencing the existing shared string).
$ws.Range("C6").Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("C3").PasteSpecial(-4163)
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("C5").PasteSpecial(-4163)

# Remote config block (rows 6-9) becomes inactive: Usar = False
# (use the stashed "False" text saved in E1 above).
$ws.Range("E1").Copy()
$ws.Range("C6").PasteSpecial(-4163)
$ws.Range("C7").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("C9").PasteSpecial(-4163)

# Clean up the scratch cell.
$ws.Range("E1").Clear()

# Update local user/password credentials from openpg/openpgpwd to odoo/odoo.
$ws.Range("B3").Value = "odoo"
$ws.Range("B4").Value = "odoo"

# Reflect where the user left the selection.
[void]$ws.Range("B5").Select()
